# Apply calibration factor update + selection change to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B3:B27 - multiply existing Position values by 0.36 (new calibration factor)
$newValues = @(7.2, 14.4, 21.6, 28.8, 36, 43.2, 50.4, 57.6, 64.8, 72, 79.2, 86.4, 93.6, 100.8, 108, 115.2, 122.4, 129.6, 136.80000000000001, 144, 151.19999999999999, 158.4, 165.6, 172.8, 180)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Update the active selection to E13
$ws.Activate()
$ws.Range("E13").Select()
